$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.369004726409912
$ws.Range("B1").Value = 4.239705085754395
$ws.Range("C1").Value = 3.662031173706055
$ws.Range("D1").Value = 1.79052460193634
$ws.Range("E1").Value = 0.9953613877296448
